$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1752863333333333
$ws.Range("H2").Value = 0.525859
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.345785
$ws.Range("N2").Value = 1.037355
$ws.Range("O2").Value = 0.1901295499320662
$ws.Range("P2").Value = 0.1901295499320662
$ws.Range("Q2").Value = 0.06061138477166667
$ws.Range("R2").Value = 0.5455024629449999
$ws.Range("S2").Value = 0.1901295499320662
$ws.Range("T2").Value = 0.1901295499320662

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1752863333333333
$ws.Range("H3").Value = 0.525859
$ws.Range("O3").Value = 0.04212778381695306
$ws.Range("P3").Value = 0.04212778381695306
$ws.Range("Q3").Value = 0.013429913001
$ws.Range("R3").Value = 0.120869217009
$ws.Range("S3").Value = 0.04212778381695306
$ws.Range("T3").Value = 0.04212778381695306

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1752863333333333
$ws.Range("H4").Value = 0.525859
$ws.Range("M4").Value = 1.396279
$ws.Range("N4").Value = 4.188836999999999
$ws.Range("O4").Value = 0.7677426662509808
$ws.Range("P4").Value = 0.7677426662509808
$ws.Range("Q4").Value = 0.2447486262203333
$ws.Range("R4").Value = 2.202737635983
$ws.Range("S4").Value = 0.7677426662509808
$ws.Range("T4").Value = 0.7677426662509808
